# Removes spell/grammar-check proofing marks (w:proofErr) that had split
# several sentences into multiple runs. Re-running Find & Replace over the
# full (already correct) visible text causes Word to re-merge the split
# runs and drop the now-orphaned w:proofErr start/end markers, matching
# the cleaned-up OOXML produced by the author's edit.

$d = $word.ActiveDocument

$enDash = [char]0x2013
$rsquo  = [char]0x2019

# Title heading: "Python " + "Youtube" + " Downloader v2.0"
$t1 = "Python Youtube Downloader v2.0"
$d.Content.Find.Execute($t1, $true, $false, $false, $false, $false, $true, 1, $false, $t1, 2) | Out-Null

# Author line: "Author: Rahul " + "Sinha"
$t2 = "Author: Rahul Sinha"
$d.Content.Find.Execute($t2, $true, $false, $false, $false, $false, $true, 1, $false, $t2, 2) | Out-Null

# Direct Video Links intro paragraph (contains "youtube" mid-sentence)
$t3 = "Cannot find a video that you wish to download? You could paste your direct link for the video and the application will try downloading it. Note: not all videos from other sources can be downloaded; all youtube links will be OK to download unless they are copyrighted."
$d.Content.Find.Execute($t3, $true, $false, $false, $false, $false, $true, 1, $false, $t3, 2) | Out-Null

# "(2) Link:" paragraph up to the hyperlink (trailing space kept before link)
$t4 = "(2) Link: This is the direct video link, for youtube this is the link like: "
$d.Content.Find.Execute($t4, $true, $false, $false, $false, $false, $true, 1, $false, $t4, 2) | Out-Null

# "Folders in the library view has following options:" (gramStart/gramEnd)
$t5 = "Folders in the library view has following options:"
$d.Content.Find.Execute($t5, $true, $false, $false, $false, $false, $true, 1, $false, $t5, 2) | Out-Null

# "Delete" bullet mentioning "infact"
$t6 = "Delete " + $enDash + " Same as above, please note for this to work the folder should be empty. If this hasn" + $rsquo + "t worked, browse to the folder using Windows and check if it is infact empty."
$d.Content.Find.Execute($t6, $true, $false, $false, $false, $false, $true, 1, $false, $t6, 2) | Out-Null
